$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin name and link for row 51 (Frax -> TheSandbox)
$ws.Range('B51').Value = 'TheSandbox'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'

# Update Price (D) and Volume(1h) (E) columns for each row
$ws.Range('D2').Value = '29.196.04'
$ws.Range('E2').Value = '  -1.18%  '
$ws.Range('D3').Value = '1.866.30'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9993'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7118'
$ws.Range('E5').Value = '  -0.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.63'
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9997'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3110'
$ws.Range('E8').Value = '  -0.39%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07673'
$ws.Range('E9').Value = '  -3.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.75'
$ws.Range('E10').Value = '  -2.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08345'
$ws.Range('E11').Value = '  +0.65%  '
$ws.Range('D12').Value = '1.868.13'
$ws.Range('E12').Value = '  -0.66%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.225'
$ws.Range('E13').Value = '  -1.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7107'
$ws.Range('E14').Value = '  -2.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.33'
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('D16').Value = '29.218.56'
$ws.Range('E16').Value = '  -1.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.929'
$ws.Range('E17').Value = '  -0.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '243.77'
$ws.Range('E18').Value = '  -1.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007816'
$ws.Range('E19').Value = '  -1.04%  '
$ws.Range('D20').Value = '2.114.92'
$ws.Range('E20').Value = '  -0.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.10'
$ws.Range('E21').Value = '  -2.19%  '
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.872'
$ws.Range('E23').Value = '  -1.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9989'
$ws.Range('E24').Value = '  -0.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1595'
$ws.Range('E25').Value = '  -1.88%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.79'
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.960'
$ws.Range('E27').Value = '  -1.42%  '
$ws.Range('E28').Value = '  +0.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.323'
$ws.Range('E29').Value = '  -2.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.499'
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.401'
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.252'
$ws.Range('E32').Value = '  +3.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05163'
$ws.Range('E33').Value = '  -2.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.8044'
$ws.Range('E34').Value = '  +10.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.914'
$ws.Range('E35').Value = '  -2.11%  '
$ws.Range('E36').Value = '  -3.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.686'
$ws.Range('E37').Value = '  +0.35%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01851'
$ws.Range('E38').Value = '  -1.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.710'
$ws.Range('E39').Value = '  -0.68%  '
$ws.Range('D40').Value = '1.165.22'
$ws.Range('E40').Value = '  -5.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.267'
$ws.Range('E41').Value = '  +0.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8968'
$ws.Range('E42').Value = '  -1.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '73.19'
$ws.Range('E43').Value = '  -2.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9991'
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '102.76'
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('D46').Value = '2.011.02'
$ws.Range('E46').Value = '  -1.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5184'
$ws.Range('E47').Value = '  -2.04%  '
$ws.Range('E48').Value = '  -0.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.343'
$ws.Range('E49').Value = '  -0.24%  '
$ws.Range('E50').Value = '  -0.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4292'
$ws.Range('E51').Value = '  -1.17%  '
